$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column AF (day 31) rows 2-6
$afValues = @{
    2 = 13791.76
    3 = 9630
    4 = 2774.9
    5 = 1228
    6 = 27424.66
}

# New values for column AG (total) rows 2-6
$agValues = @{
    2 = 309737.14
    3 = 141201.81
    4 = 90788.8
    5 = 73082.28999999999
    6 = 614810.04
}

foreach ($row in 2..6) {
    $ws.Range("AF$row").Value = $afValues[$row]
    $ws.Range("AG$row").Value = $agValues[$row]
}
